# Applies the Dec 30 2023 cryptos data refresh described in the commit.
# Every target cell is forced to Text format before the write (many of the
# "Price" values look like plain numbers, e.g. "39.50", and Excel would
# otherwise silently coerce them to the number 39.5) and then the style is
# reset back to Normal so no stray number-format/quote-prefix style lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.291.46"
Set-TextValue "E2" "  +0.75%  "
Set-TextValue "D3" "2.300.89"
Set-TextValue "E3" "  -0.05%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "316.43"
Set-TextValue "E5" "  +1.11%  "
Set-TextValue "D6" "102.66"
Set-TextValue "E6" "  -1.80%  "
Set-TextValue "D7" "0.623"
Set-TextValue "E7" "  +0.39%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "D9" "0.605"
Set-TextValue "E9" "  -0.31%  "
Set-TextValue "D10" "39.50"
Set-TextValue "E10" "  -1.98%  "
Set-TextValue "D11" "0.0908"
Set-TextValue "E11" "  -0.50%  "
Set-TextValue "D12" "8.38"
Set-TextValue "E12" "  +1.29%  "
Set-TextValue "D13" "0.106"
Set-TextValue "E13" "  +0.44%  "
Set-TextValue "D14" "0.961"
Set-TextValue "E14" "  -1.25%  "
Set-TextValue "D15" "15.20"
Set-TextValue "E15" "  -2.39%  "
Set-TextValue "D16" "2.649.67"
Set-TextValue "E16" "  +0.17%  "
Set-TextValue "D17" "2.291.26"
Set-TextValue "E17" "  -0.18%  "
Set-TextValue "D18" "42.393.08"
Set-TextValue "E18" "  +0.82%  "
Set-TextValue "D19" "7.44"
Set-TextValue "E19" "  -2.31%  "
Set-TextValue "E20" "  +0.75%  "
Set-TextValue "D21" "73.42"
Set-TextValue "E21" "  -1.50%  "
Set-TextValue "D22" "3.54"
Set-TextValue "E22" "  +2.56%  "
Set-TextValue "D23" "276.45"
Set-TextValue "E23" "  +7.35%  "
Set-TextValue "D24" "11.35"
Set-TextValue "E24" "  +21.87%  "
Set-TextValue "D25" "2.26"
Set-TextValue "E25" "  -1.75%  "
Set-TextValue "E27" "  -1.22%  "
Set-TextValue "E28" "  +3.24%  "
Set-TextValue "D29" "22.73"
Set-TextValue "E29" "  -0.38%  "
Set-TextValue "D30" "37.49"
Set-TextValue "E30" "  +5.02%  "
Set-TextValue "D31" "165.69"
Set-TextValue "E31" "  -0.64%  "
Set-TextValue "D32" "0.0873"
Set-TextValue "E32" "  -2.60%  "
Set-TextValue "D33" "5.86"
Set-TextValue "E33" "  +0.87%  "
Set-TextValue "E34" "  +4.66%  "
Set-TextValue "E35" "  -0.43%  "
Set-TextValue "E36" "  -10.46%  "
Set-TextValue "E37" "  +3.62%  "
Set-TextValue "D38" "4.57"
Set-TextValue "E38" "  +0.23%  "
Set-TextValue "D39" "3.69"
Set-TextValue "E39" "  +1.52%  "
Set-TextValue "D40" "2.77"
Set-TextValue "E40" "  -0.21%  "
Set-TextValue "D41" "1.50"
Set-TextValue "E41" "  +3.24%  "
Set-TextValue "D42" "69.67"
Set-TextValue "E42" "  -3.03%  "
Set-TextValue "B43" "BitcoinSV"
Set-TextValue "C43" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D43" "94.83"
Set-TextValue "E43" "  -3.69%  "
Set-TextValue "B44" "Algorand"
Set-TextValue "C44" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D44" "0.226"
Set-TextValue "E44" "  -0.45%  "
Set-TextValue "E45" "  +0.13%  "
Set-TextValue "D46" "80.92"
Set-TextValue "E46" "  +8.45%  "
Set-TextValue "D47" "12.03"
Set-TextValue "E47" "  -2.27%  "
Set-TextValue "D48" "113.01"
Set-TextValue "E48" "  +0.51%  "
Set-TextValue "D49" "8.96"
Set-TextValue "E49" "  -0.90%  "
Set-TextValue "E50" "  -1.47%  "
Set-TextValue "D51" "1.588.51"
Set-TextValue "E51" "  +1.26%  "

Write-Output "Updated cryptos list: applied $(91) cell changes (rows 2-51)."
